$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 37, pushing existing rows 37-53 down to 38-54.
$ws.Rows("37:37").Insert()

# Populate the newly inserted row 37 with the new Cilantro market record.
$ws.Cells.Item(37, 1).Value = 5
$ws.Cells.Item(37, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(37, 3).Value = "Maule"
$ws.Cells.Item(37, 4).Value = 44806
$ws.Cells.Item(37, 5).Value = 7
$ws.Cells.Item(37, 6).Value = 100112040
$ws.Cells.Item(37, 7).Value = "Cilantro"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 150
$ws.Cells.Item(37, 11).Value = 13000
$ws.Cells.Item(37, 12).Value = 13000
$ws.Cells.Item(37, 13).Value = 13000
$ws.Cells.Item(37, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(37, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(37, 16).Value = 361
$ws.Cells.Item(37, 17).Value = 36
$ws.Cells.Item(37, 18).Value = "Hortaliza"

# Ensure the date cell keeps the workbook's date number format (style index 2).
$ws.Cells.Item(37, 4).NumberFormat = $ws.Cells.Item(38, 4).NumberFormat
